# Adjust Investment Summary table column widths for better formatting
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$tbl = $s.Shapes.Item(3).Table

# New column widths in EMU, converted to points (1 pt = 12700 EMU)
$widthsEmu = @(1742186, 1045311, 2003514, 1132421, 871093, 871093, 1045311)

for ($i = 0; $i -lt $widthsEmu.Length; $i++) {
    $tbl.Columns.Item($i + 1).Width = $widthsEmu[$i] / 12700
}
